$d = $word.ActiveDocument

# The template's signature/date placeholder line looks like:
#   "…………, ……......……. 20..."
# and needs to become the bracket-token form used by the live form:
#   "#2#, #3#"   (i.e. run "#2#", run ", ", run "#3", <bookmark>, run "#")
#
# wdReplaceNone = 0, wdFindContinue = 1

# --- Locate & rewrite the first chunk "…………, ……" -> "#2#" --------------
$rFirst = $d.Content
$rFirst.Find.Execute("…………, ……", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rFirst.Text = "#2#"

# --- Locate the remainder of the old placeholder ("......…….<space>") and
#     delete it outright; it sits immediately after the text we just set. --
$afterFirst = $d.Range($rFirst.End, $d.Content.End)
$rMid = $afterFirst
$rMid.Find.Execute("......…….", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# extend the match by the trailing single space run that follows it
$rMid = $d.Range($rMid.Start, $rMid.End + 1)
$insPos = $rMid.Start
$rMid.Text = ""

# --- Insert the two new runs ", " and "#3" in the gap ---------------------
# (A collapsed-range insert always glues onto the immediately preceding run,
#  so splice both pieces of text in first, then toggle Bold off/on over each
#  piece to force Word to keep them as their own runs instead of coalescing
#  back into one.)
$rInsA = $d.Range($insPos, $insPos)
$rInsA.InsertAfter(", ")
$rInsB = $d.Range($insPos + 2, $insPos + 2)
$rInsB.InsertAfter("#3")

$rSplitRight = $d.Range($insPos + 2, $insPos + 4)
$rSplitRight.Bold = 1
$rSplitRight.Bold = 0
$rSplitLeft = $d.Range($insPos, $insPos + 2)
$rSplitLeft.Bold = 1
$rSplitLeft.Bold = 0

# --- Final chunk "20..." -> "#" --------------------------------------------
$d.Content.Find.Execute("20...", $false, $false, $false, $false, $false, $true, 1, $false, "#", 2) | Out-Null

Write-Host "edit applied"
